# Accept project path as argument and change the column Severity to Action.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workflow")

# Rename the column header from "Severity" to "Action"
$ws.Range("E1").Value = "Action"

# Update the per-row action values (previously Severity levels Low/Medium/High)
$ws.Range("E2").Value = "Fix"
$ws.Range("E3").Value = "Fix"
$ws.Range("E4").Value = "Fix"
$ws.Range("E5").Value = "Fix"
$ws.Range("E6").Value = "Fix"
$ws.Range("E7").Value = "Fix"
$ws.Range("E8").Value = "Double check"
$ws.Range("E9").Value = "Double check"
$ws.Range("E10").Value = "Double check"
$ws.Range("E11").Value = "Double check"
$ws.Range("E12").Value = "Double check"

# Update the active selection to match the post-edit UI state
$ws.Range("D2").Select() | Out-Null
